# Update vehicle specifications from EV database, Mercedes Benz, Volkswagen
#
# The ICEV "Use phase" row (previously row 9, raw t CO2 values) is pushed
# down to a brand-new row 11, and row 9 is repurposed to show the
# equivalent emissions expressed as g CO2/km (computed from the values
# now living in row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Insert a new row 11 and move the old "Use phase" / "t CO2" figures
#    (formerly D9:G9) down there, keeping their original 0.0 number format.
$ws.Rows.Item(11).Insert()

$ws.Range("B11").Value = "Use phase"
$ws.Range("C11").Value = "t CO2"
$ws.Range("D11").Value = 23.736457765667577
$ws.Range("E11").Value = 30.98688033233573
$ws.Range("F11").Value = 35.047812067271181
$ws.Range("G11").Value = 45.788684113548037
$ws.Range("D11:G11").NumberFormat = "0.0"

# 2) Row 9 becomes the "g CO2/km" row, computed from row 11's t CO2 values
#    (t CO2 per vehicle over a 180,000 km lifetime -> g CO2/km).
$ws.Range("C9").Value = "g CO2/km"

$ws.Range("D9").Formula = "=D11*1000000/180000"
$ws.Range("E9").Formula = "=E11*1000000/180000"
$ws.Range("F9").Formula = "=F11*1000000/180000"
$ws.Range("G9").Formula = "=G11*1000000/180000"

$ws.Range("D9:G9").NumberFormat = "0.0"
$ws.Range("D9:G9").HorizontalAlignment = -4152

# 3) Columns D:G get a uniform width instead of per-column bestFit widths.
$ws.Range("D1:G1").ColumnWidth = 9.65

$ws.Calculate()

# 4) Match the author's last on-screen selection.
[void]$ws.Range("D13:J15").Select()
